# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G on Sheet1 is labeled "K" (strikeouts) in row 1. The previous
# export populated it from a "Strike#" style metric; this regenerates the
# column from the corrected K values (pitch-by-pitch K count) for every
# data row (rows 2-73).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1;  3  = 1;  4  = 1;  5  = 0;  6  = 2;  7  = 0;  8  = 0;  9  = 0;
    10 = 1;  11 = 3;  12 = 1;  13 = 1;  14 = 2;  15 = 3;  16 = 2;  17 = 2;
    18 = 0;  19 = 1;  20 = 0;  21 = 1;  22 = 2;  23 = 1;  24 = 4;  25 = 1;
    26 = 0;  27 = 1;  28 = 0;  29 = 0;  30 = 0;  31 = 1;  32 = 1;  33 = 0;
    34 = 1;  35 = 1;  36 = 1;  37 = 1;  38 = 2;  39 = 0;  40 = 3;  41 = 0;
    42 = 0;  43 = 0;  44 = 1;  45 = 2;  46 = 0;  47 = 1;  48 = 0;  49 = 1;
    50 = 0;  51 = 0;  52 = 0;  53 = 2;  54 = 1;  55 = 1;  56 = 3;  57 = 1;
    58 = 1;  59 = 1;  60 = 1;  61 = 1;  62 = 1;  63 = 1;  64 = 0;  65 = 0;
    66 = 1;  67 = 2;  68 = 0;  69 = 0;  70 = 0;  71 = 0;  72 = 0;  73 = 1;
}

# Column G is the 7th column (A=1 ... J=10).
$col = 7

foreach ($row in ($kValues.Keys | Sort-Object)) {
    $ws.Cells.Item($row, $col).Value = $kValues[$row]
}
